$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Massachusetts)
$ws.Range("B3").Value = 44023
$ws.Range("C3").Value = 111398
$ws.Range("D3").Value = 8310
$ws.Range("E3").Value = 10477
$ws.Range("F3").Value = 684
$ws.Range("G3").Value = 9.41

# Row 5 (Iowa) -- B5 already 44023, unchanged
$ws.Range("C5").Value = 34647
$ws.Range("E5").Value = 3010

# Row 12 (Wisconsin)
$ws.Range("B12").Value = 44023
$ws.Range("C12").Value = 35679
$ws.Range("D12").Value = 821
$ws.Range("E12").Value = 5984
$ws.Range("G12").Value = 18.6
$ws.Range("H12").Value = 23.73
$ws.Range("K12").Value = 32164
$ws.Range("L12").Value = 809

# Row 14 (Tennessee)
$ws.Range("B14").Value = 44023
$ws.Range("C14").Value = 61006
$ws.Range("D14").Value = 738
$ws.Range("E14").Value = 12323
$ws.Range("F14").Value = 261
$ws.Range("G14").Value = 20.2
$ws.Range("H14").Value = 35.37

# Row 16 (Utah) -- C/D/E are stored as text (inlineStr) values that look numeric
$ws.Range("B16").Value = 44023
$ws.Range("C16").Value = "'28855"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'212"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'731"
$ws.Range("E16").Style = "Normal"

# Row 21 (Alaska)
$ws.Range("B21").Value = 44023
$ws.Range("C21").Value = 1385
$ws.Range("E21").Value = 34
$ws.Range("G21").Value = 1.43
$ws.Range("K21").Value = 2376

# Row 34 (Georgia)
$ws.Range("B34").Value = 44023
$ws.Range("C34").Value = 114401
$ws.Range("D34").Value = 2996
$ws.Range("E34").Value = 31052
$ws.Range("F34").Value = 1397
$ws.Range("G34").Value = 27.14
$ws.Range("H34").Value = 46.63

# Row 41 (Michigan)
$ws.Range("B41").Value = 44023
$ws.Range("C41").Value = 68857
$ws.Range("D41").Value = 5983
$ws.Range("E41").Value = 20486
$ws.Range("F41").Value = 2390
$ws.Range("G41").Value = 29.75
$ws.Range("H41").Value = 39.95
